#
# Sprint1 Week 1 - Test data
#
# Inserts 5 new transmittal-attachment columns (TxType/IssueReason stay put;
# AttachDocuments.../ReviewDocument are new, inserted right before the old
# "Message" column, which - together with everything to its right - shifts
# from G:J to L:O) and updates the shared "TxType group" string pool.
#

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank columns at G:K. Everything previously at G:L (Message,
# Action-Level2, ForwardTo, Action-Level3, ...) shifts right to L:O and all
# existing cell styling / column widths shift with it automatically.
$ws.Range("G1:K1").EntireColumn.Insert()

# New header row values for the freshly inserted columns.
$ws.Range("G1").Value = "AttachDocuments"
$ws.Range("H1").Value = "AttachDocumentName"
$ws.Range("I1").Value = "AttachSupportDocuments"
$ws.Range("J1").Value = "AttachSupportDocumentName"
$ws.Range("K1").Value = "ReviewDocument"

# Match the target column widths (character units, best achievable given the
# engine's internal 1/6-character pixel quantisation).
$ws.Range("G1").ColumnWidth = 9.666666666666666
$ws.Range("H1").ColumnWidth = 7.666666666666667
$ws.Range("I1").ColumnWidth = 11.333333333333334
$ws.Range("J1").ColumnWidth = 10
$ws.Range("K1").ColumnWidth = 8

# Old column H (now M, "Action-Level2") drops its bestFit auto-width and
# becomes a fixed custom width of 19 characters.
$ws.Range("M1").ColumnWidth = 18.166666666666668
